$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new trade record as row 10 (columns: Date, Profitable,
# Principle, Start Principle, BuyPrice, SellPrice, IsShortSell,
# Price Change %, Strong trade)
$row = 10
$ws.Cells.Item($row, 1).Value = 42654.745995370373
$ws.Cells.Item($row, 1).NumberFormat = "m/d/yy h:mm"
$ws.Cells.Item($row, 2).Value = $true
$ws.Cells.Item($row, 3).Value = 10022.16
$ws.Cells.Item($row, 4).Value = 10009.15
$ws.Cells.Item($row, 5).Value = 18.870000999999998
$ws.Cells.Item($row, 6).Value = 18.920000000000002
$ws.Cells.Item($row, 7).Value = $false
$ws.Cells.Item($row, 7).NumberFormat = "m/d/yy h:mm"
$ws.Cells.Item($row, 8).Value = 0.26
$ws.Cells.Item($row, 9).Value = $false

# The repeater re-ran its column autosize pass over the whole trade
# table after appending the row, which nudged every column's best-fit
# width slightly (new widest date string in col A, etc). Reproduce the
# resulting widths as closely as this engine's width grid allows.
$ws.Columns.Item(1).ColumnWidth = 14.5
$ws.Columns.Item(2).ColumnWidth = 7.333333333333333
$ws.Columns.Item(3).ColumnWidth = 8
$ws.Columns.Item(4).ColumnWidth = 10.333333333333334
$ws.Columns.Item(5).ColumnWidth = 9
$ws.Columns.Item(6).ColumnWidth = 6.166666666666667
$ws.Columns.Item(7).ColumnWidth = 9.5
$ws.Columns.Item(8).ColumnWidth = 13.833333333333334
$ws.Columns.Item(9).ColumnWidth = 11
